$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.757.68"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "2.832.45"
$ws.Range("E3").Value = "  +2.58%  "
$ws.Range("D5").Value = "'352.43"
$ws.Range("E5").Value = "  +5.91%  "
$ws.Range("D6").Value = "'113.71"
$ws.Range("E6").Value = "  -2.44%  "
$ws.Range("D7").Value = "'0.557"
$ws.Range("E7").Value = "  +3.34%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.599"
$ws.Range("E9").Value = "  +3.84%  "
$ws.Range("E10").Value = "  -0.50%  "
$ws.Range("D11").Value = "'0.0851"
$ws.Range("E11").Value = "  -0.92%  "
$ws.Range("E12").Value = "  -1.00%  "
$ws.Range("E13").Value = "  +1.40%  "
$ws.Range("E14").Value = "  +0.89%  "
$ws.Range("D15").Value = "3.262.87"
$ws.Range("D16").Value = "2.822.49"
$ws.Range("E16").Value = "  +1.70%  "
$ws.Range("D17").Value = "'0.894"
$ws.Range("E17").Value = "  +0.65%  "
$ws.Range("D18").Value = "51.656.37"
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("D19").Value = "'7.40"
$ws.Range("E19").Value = "  +7.73%  "
$ws.Range("D20").Value = "'3.14"
$ws.Range("E20").Value = "  -2.67%  "
$ws.Range("D21").Value = "'13.45"
$ws.Range("E21").Value = "  -0.36%  "
$ws.Range("E22").Value = "  +1.98%  "
$ws.Range("D23").Value = "'270.33"
$ws.Range("E23").Value = "  -2.87%  "
$ws.Range("D24").Value = "'69.67"
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("E25").Value = "  +3.60%  "
$ws.Range("D26").Value = "'26.72"
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("D28").Value = "'10.28"
$ws.Range("E28").Value = "  +1.11%  "
$ws.Range("E29").Value = "  +1.23%  "
$ws.Range("E30").Value = "  -1.50%  "
$ws.Range("D31").Value = "'50.67"
$ws.Range("E31").Value = "  +1.32%  "
$ws.Range("E32").Value = "  -3.77%  "
$ws.Range("D33").Value = "'0.0449"
$ws.Range("E33").Value = "  +27.68%  "
$ws.Range("D34").Value = "'5.80"
$ws.Range("E34").Value = "  +4.22%  "
$ws.Range("D35").Value = "'0.0825"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("D37").Value = "'2.08"
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "'3.22"
$ws.Range("E38").Value = "  -0.75%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "'4.87"
$ws.Range("E39").Value = "  -2.99%  "
$ws.Range("D40").Value = "'18.12"
$ws.Range("E40").Value = "  -4.59%  "
$ws.Range("D41").Value = "'23.91"
$ws.Range("E41").Value = "  +3.30%  "
$ws.Range("D42").Value = "'2.55"
$ws.Range("E42").Value = "  +4.72%  "
$ws.Range("E43").Value = "  +1.00%  "
$ws.Range("D44").Value = "'125.44"
$ws.Range("E44").Value = "  -1.13%  "
$ws.Range("E45").Value = "  +0.34%  "
$ws.Range("D46").Value = "2.082.48"
$ws.Range("E46").Value = "  -0.31%  "
$ws.Range("E47").Value = "  +0.54%  "
$ws.Range("E48").Value = "  +3.59%  "
$ws.Range("E49").Value = "  +2.68%  "
$ws.Range("D50").Value = "'0.942"
$ws.Range("E50").Value = "  +7.75%  "
$ws.Range("D51").Value = "'60.56"
$ws.Range("E51").Value = "  +1.00%  "
